$wb = $excel.ActiveWorkbook

# ---- GLOBAL RESULTS sheet ----
$ws = $wb.Worksheets.Item("GLOBAL RESULTS")

$ws.Range("C2").Value  = 60.540708391007826
$ws.Range("C3").Value  = 12.578250317120148
$ws.Range("C4").Value  = 24.183035275777005
$ws.Range("C5").Value  = 0.5567940147707489

$ws.Range("C7").Value  = 43.608243744263845
$ws.Range("C8").Value  = 12.188394583002776
$ws.Range("C9").Value  = 26.8984357935812
$ws.Range("C10").Value = 0.6193138241651133

$ws.Range("C12").Value = 43.608243744263845
$ws.Range("C13").Value = 12.188394583002776
$ws.Range("C14").Value = 26.8984357935812
$ws.Range("C15").Value = 0.6193138241651133

$ws.Range("C17").Value = 42.16856219521241
$ws.Range("C18").Value = 12.155247128163197
$ws.Range("C19").Value = 17.525864882468294
$ws.Range("C20").Value = 0.40351827464824574

$ws.Range("C22").Value = 41.66724914526735
$ws.Range("C23").Value = 12.143704817508112
$ws.Range("C24").Value = 24.712347950739662
$ws.Range("C25").Value = 0.5689809932042091

$ws.Range("C27").Value = 29.55032474357941
$ws.Range("C28").Value = 56.77536475163592

# ---- LANDING GEARS sheet ----
$ws2 = $wb.Worksheets.Item("LANDING GEARS")
$ws2.Range("C2").Value = 12.321629305144317
